# Generate Report for Handback
#
# The "9b4bbf56-..." handback row failed because the returned file name
# did not match the original handoff file name. Reflect that failure in
# the status text (Overview + per-locale sheets) and record the detailed
# error message in the "Error Detail" column for both the zh-cn and
# de-de locale sheets, widening that column so the message is readable.

$wb = $excel.ActiveWorkbook

$failedStatus = "Handback transform failed"

# --- Update every cell that shares the old "Ready for handoff" string ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $failedStatus
$wsOverview.Range("F3").Value = $failedStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $failedStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $failedStatus

# --- Fill in the Error Detail (column P) message for the failed handback ---
$wsZhCn.Range("P3").Value = "Handback file name: diepnwdg.14x is different with handoff file name: 9b4bbf56-bfe8-4fc7-8587-12838d830880.eb6fc0ae24fe23334a74a131e14f21a2698b8ff5.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: diepnwdg.14x is different with handoff file name: 9b4bbf56-bfe8-4fc7-8587-12838d830880.eb6fc0ae24fe23334a74a131e14f21a2698b8ff5.de-de."

# --- Widen the Error Detail column so the message is visible ---
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
